$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4068093373250292
$ws.Range("D2").Value = 0.009351139749049109
$ws.Range("E2").Value = 0.1720958756594229
$ws.Range("F2").Value = 0.9159663923201435
$ws.Range("G2").Value = 0.7963027762735919
$ws.Range("H2").Value = 0.7612044533489666
$ws.Range("L2").Value = 0.15040379545691
$ws.Range("O2").Value = 3.12815358371347
$ws.Range("C3").Value = 0.3975481905092693
$ws.Range("D3").Value = 0.009435166500940673
$ws.Range("E3").Value = 0.1677685677011951
$ws.Range("F3").Value = 0.8772088377580189
$ws.Range("G3").Value = 0.7538445135139114
$ws.Range("H3").Value = 0.7461335867586456
$ws.Range("L3").Value = 0.1462578624851787
$ws.Range("O3").Value = 3.006157524454125
$ws.Range("C4").Value = 0.3920966218257433
$ws.Range("D4").Value = 0.00948832938111499
$ws.Range("E4").Value = 0.1652137714663446
$ws.Range("F4").Value = 0.8539758345155519
$ws.Range("G4").Value = 0.7282703734057634
$ws.Range("H4").Value = 0.7373056249113006
$ws.Range("L4").Value = 0.1438038777468833
$ws.Range("O4").Value = 2.933175414339644
$ws.Range("C5").Value = 0.3899339996491733
$ws.Range("D5").Value = 0.009510388064638331
$ws.Range("E5").Value = 0.1641983007239354
$ws.Range("F5").Value = 0.8446496412524453
$ws.Range("G5").Value = 0.7179726945568348
$ws.Range("H5").Value = 0.7338150516630435
$ws.Range("L5").Value = 0.1428268243803146
$ws.Range("O5").Value = 2.90391688686924
$ws.Range("C6").Value = 0.3895784543456102
$ws.Range("D6").Value = 0.009514074699972896
$ws.Range("E6").Value = 0.1640312286105896
$ws.Range("F6").Value = 0.8431095637789809
$ws.Range("G6").Value = 0.7162702446704827
$ws.Range("H6").Value = 0.7332418969721743
$ws.Range("L6").Value = 0.1426659702311852
$ws.Range("O6").Value = 2.899087614875839
$ws.Range("C7").Value = 0.392067217457992
$ws.Range("D7").Value = 0.009488625276098617
$ws.Range("E7").Value = 0.165199972786958
$ws.Range("F7").Value = 0.8538494861101498
$ws.Range("G7").Value = 0.7281309940444487
$ws.Range("H7").Value = 0.7372581171931074
$ws.Range("L7").Value = 0.1437906079781754
$ws.Range("O7").Value = 2.932778872918163
$ws.Range("C8").Value = 0.4035672883939014
$ws.Range("D8").Value = 0.009379786018830716
$ws.Range("E8").Value = 0.1705825686378546
$ws.Range("F8").Value = 0.9024853816899849
$ws.Range("G8").Value = 0.7815599795246442
$ws.Range("H8").Value = 0.7559195968969448
$ws.Range("L8").Value = 0.1489552107994783
$ws.Range("O8").Value = 3.085689059124093
$ws.Range("C9").Value = 0.4279886984085124
$ws.Range("D9").Value = 0.009178828835648689
$ws.Range("E9").Value = 0.1819524292491792
$ws.Range("F9").Value = 1.002364138748504
$ws.Range("G9").Value = 0.8902974908077681
$ws.Range("H9").Value = 0.7959006314800376
$ws.Range("L9").Value = 0.1598142602096431
$ws.Range("O9").Value = 3.400899602077232
$ws.Range("C10").Value = 0.4470827431894691
$ws.Range("D10").Value = 0.009038814222902536
$ws.Range("E10").Value = 0.1908088464605058
$ws.Range("F10").Value = 1.078536330038418
$ws.Range("G10").Value = 0.9726574500497236
$ws.Range("H10").Value = 0.8273555155267047
$ws.Range("L10").Value = 0.1682450686695205
$ws.Range("O10").Value = 3.641995857184725
$ws.Range("C11").Value = 0.4560218970189283
$ws.Range("D11").Value = 0.008976781069067741
$ws.Range("E11").Value = 0.194948531687821
$ws.Range("F11").Value = 1.113805546471013
$ws.Range("G11").Value = 1.010673884672855
$ws.Range("H11").Value = 0.8421208982347821
$ws.Range("L11").Value = 0.1721803059861031
$ws.Range("O11").Value = 3.753776474961853
$ws.Range("C12").Value = 0.4594434993398124
$ws.Range("D12").Value = 0.00895353041766489
$ws.Range("E12").Value = 0.1965321595870222
$ws.Range("F12").Value = 1.127250626168475
$ws.Range("G12").Value = 1.025149684592009
$ws.Range("H12").Value = 0.8477780284605672
$ws.Range("L12").Value = 0.1736849660905335
$ws.Range("O12").Value = 3.796409758596155
$ws.Range("C13").Value = 0.4587049688618947
$ws.Range("D13").Value = 0.008958527168961616
$ws.Range("E13").Value = 0.1961903832866838
$ws.Range("F13").Value = 1.124351003731732
$ws.Range("G13").Value = 1.022028502582913
$ws.Range("H13").Value = 0.8465567350012009
$ws.Range("L13").Value = 0.1733602655907163
$ws.Range("O13").Value = 3.787214352518163
$ws.Range("C14").Value = 0.4563026612454735
$ws.Range("D14").Value = 0.008974863407454237
$ws.Range("E14").Value = 0.1950784962148049
$ws.Range("F14").Value = 1.11490988676033
$ws.Range("G14").Value = 1.011863214402723
$ws.Range("H14").Value = 0.8425849940172156
$ws.Range("L14").Value = 0.1723038047452974
$ws.Range("O14").Value = 3.757277830442035
$ws.Range("C15").Value = 0.4548359408504723
$ws.Range("D15").Value = 0.008984901113034294
$ws.Range("E15").Value = 0.1943995217469379
$ws.Range("F15").Value = 1.109138586840345
$ws.Range("G15").Value = 1.005647093545633
$ws.Range("H15").Value = 0.8401607622488427
$ws.Range("L15").Value = 0.1716585791814822
$ws.Range("O15").Value = 3.738980538580222
$ws.Range("C16").Value = 0.446503654289728
$ws.Range("D16").Value = 0.009042901745627052
$ws.Range("E16").Value = 0.1905405456508973
$ws.Range("F16").Value = 1.076243885124114
$ws.Range("G16").Value = 0.9701841184887883
$ws.Range("H16").Value = 0.8263997614475045
$ws.Range("L16").Value = 0.1679899118500146
$ws.Range("O16").Value = 3.634733231000098
$ws.Range("C17").Value = 0.4414570045306618
$ws.Range("D17").Value = 0.009078909304715221
$ws.Range("E17").Value = 0.1882016369687989
$ws.Range("F17").Value = 1.056222723677024
$ws.Range("G17").Value = 0.9485701779534281
$ws.Range("H17").Value = 0.8180748479653914
$ws.Range("L17").Value = 0.1657649834149737
$ws.Range("O17").Value = 3.571321180802784
$ws.Range("C18").Value = 0.4385781272709721
$ws.Range("D18").Value = 0.009099776085649491
$ws.Range("E18").Value = 0.1868667881919706
$ws.Range("F18").Value = 1.044765229239204
$ws.Range("G18").Value = 0.9361901856157431
$ws.Range("H18").Value = 0.8133295258420503
$ws.Range("L18").Value = 0.1644946724680949
$ws.Range("O18").Value = 3.535046248280196
$ws.Range("C19").Value = 0.4376074763875692
$ws.Range("D19").Value = 0.009106868009119751
$ws.Range("E19").Value = 0.1864166202737749
$ws.Range("F19").Value = 1.040895890022554
$ws.Range("G19").Value = 0.9320074057544048
$ws.Range("H19").Value = 0.811730212508138
$ws.Range("L19").Value = 0.1640661803442924
$ws.Range("O19").Value = 3.522798133886397
$ws.Range("C20").Value = 0.4419917625536698
$ws.Range("D20").Value = 0.00907506006614911
$ws.Range("E20").Value = 0.1884495378012971
$ws.Range("F20").Value = 1.058347988721053
$ws.Range("G20").Value = 0.9508656557133293
$ws.Range("H20").Value = 0.8189566030620483
$ws.Range("L20").Value = 0.1660008564074218
$ws.Range("O20").Value = 3.578050996948036
$ws.Range("C21").Value = 0.4570072845094728
$ws.Range("D21").Value = 0.008970058534938818
$ws.Range("E21").Value = 0.1954046490250718
$ws.Range("F21").Value = 1.11768054147862
$ws.Range("G21").Value = 1.014846835011895
$ws.Range("H21").Value = 0.8437498034405735
$ws.Range("L21").Value = 0.1726137194532242
$ws.Range("O21").Value = 3.766062637682921
$ws.Range("C22").Value = 0.4670338348652479
$ws.Range("D22").Value = 0.008902833278146671
$ws.Range("E22").Value = 0.2000436238092291
$ws.Range("F22").Value = 1.156979085665597
$ws.Range("G22").Value = 1.057127679719855
$ws.Range("H22").Value = 0.8603372093284065
$ws.Range("L22").Value = 0.1770199983058944
$ws.Range("O22").Value = 3.890714148656627
$ws.Range("C23").Value = 0.4616629389621778
$ws.Range("D23").Value = 0.008938584190553822
$ws.Range("E23").Value = 0.1975591428688688
$ws.Range("F23").Value = 1.135956844762489
$ws.Range("G23").Value = 1.034518794695117
$ws.Range("H23").Value = 0.8514490402922377
$ws.Range("L23").Value = 0.1746605331655644
$ws.Range("O23").Value = 3.824022297306726
$ws.Range("C24").Value = 0.4417499283512711
$ws.Range("D24").Value = 0.009076799791421841
$ws.Range("E24").Value = 0.1883374312511137
$ws.Range("F24").Value = 1.057386991536589
$ws.Range("G24").Value = 0.9498277265160198
$ws.Range("H24").Value = 0.8185578345618012
$ws.Range("L24").Value = 0.1658941907309242
$ws.Range("O24").Value = 3.57500788184376
$ws.Range("C25").Value = 0.4211806766276709
$ws.Range("D25").Value = 0.009231855463133654
$ws.Range("E25").Value = 0.1787887158460251
$ws.Range("F25").Value = 0.9748574239048509
$ws.Range("G25").Value = 0.8604512084784517
$ws.Range("H25").Value = 0.7338150516630435
$ws.Range("L25").Value = 0.1567976055937805
$ws.Range("O25").Value = 3.313967518965114
